$wb = $excel.ActiveWorkbook

# --- Proximity sheet: append two rows (ENTER/EXIT) ---
$ws = $wb.Worksheets.Item("Proximity")

$ws.Cells.Item(12, 1).Value = "'2026-02-01"
$ws.Cells.Item(12, 2).Value = "15:59:51"
$ws.Cells.Item(12, 3).Value = "15:00"
$ws.Cells.Item(12, 4).Value = "Living Room Main Door"
$ws.Cells.Item(12, 5).Value = "ENTER"
$ws.Cells.Item(12, 6).Value = "User ENTERED Living Room Main Door"

$ws.Cells.Item(13, 1).Value = "'2026-02-01"
$ws.Cells.Item(13, 2).Value = "16:00:10"
$ws.Cells.Item(13, 3).Value = "16:00"
$ws.Cells.Item(13, 4).Value = "Living Room Main Door"
$ws.Cells.Item(13, 5).Value = "EXIT"
$ws.Cells.Item(13, 6).Value = "User EXITED Living Room Main Door"

# --- mmWave sheet: append six PRESENCE_DETECTED rows ---
$ws = $wb.Worksheets.Item("mmWave")

$ws.Cells.Item(48, 1).Value = "'2026-02-01"
$ws.Cells.Item(48, 2).Value = "15:59:51"
$ws.Cells.Item(48, 3).Value = "15:00"
$ws.Cells.Item(48, 4).Value = "Living Room"
$ws.Cells.Item(48, 5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(48, 6).Value = "Active"

$ws.Cells.Item(49, 1).Value = "'2026-02-01"
$ws.Cells.Item(49, 2).Value = "15:59:57"
$ws.Cells.Item(49, 3).Value = "15:00"
$ws.Cells.Item(49, 4).Value = "Living Room"
$ws.Cells.Item(49, 5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(49, 6).Value = "Active"

$ws.Cells.Item(50, 1).Value = "'2026-02-01"
$ws.Cells.Item(50, 2).Value = "16:00:07"
$ws.Cells.Item(50, 3).Value = "16:00"
$ws.Cells.Item(50, 4).Value = "Living Room"
$ws.Cells.Item(50, 5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(50, 6).Value = "Active"

$ws.Cells.Item(51, 1).Value = "'2026-02-01"
$ws.Cells.Item(51, 2).Value = "16:00:18"
$ws.Cells.Item(51, 3).Value = "16:00"
$ws.Cells.Item(51, 4).Value = "Living Room"
$ws.Cells.Item(51, 5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(51, 6).Value = "Active"

$ws.Cells.Item(52, 1).Value = "'2026-02-01"
$ws.Cells.Item(52, 2).Value = "16:00:28"
$ws.Cells.Item(52, 3).Value = "16:00"
$ws.Cells.Item(52, 4).Value = "Living Room"
$ws.Cells.Item(52, 5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(52, 6).Value = "Active"

$ws.Cells.Item(53, 1).Value = "'2026-02-01"
$ws.Cells.Item(53, 2).Value = "16:00:39"
$ws.Cells.Item(53, 3).Value = "16:00"
$ws.Cells.Item(53, 4).Value = "Living Room"
$ws.Cells.Item(53, 5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(53, 6).Value = "Active"

# --- Camera sheet: append three "Image Captured" rows ---
$ws = $wb.Worksheets.Item("Camera")

$ws.Cells.Item(10, 1).Value = "'2026-02-01"
$ws.Cells.Item(10, 2).Value = "15:59:49"
$ws.Cells.Item(10, 3).Value = "15:00"
$ws.Cells.Item(10, 4).Value = "Living Room Main Door"
$ws.Cells.Item(10, 5).Value = "Image Captured"
$ws.Cells.Item(10, 6).Value = "Active"

$ws.Cells.Item(11, 1).Value = "'2026-02-01"
$ws.Cells.Item(11, 2).Value = "15:59:51"
$ws.Cells.Item(11, 3).Value = "15:00"
$ws.Cells.Item(11, 4).Value = "Living Room Main Door"
$ws.Cells.Item(11, 5).Value = "Image Captured"
$ws.Cells.Item(11, 6).Value = "Active"

$ws.Cells.Item(12, 1).Value = "'2026-02-01"
$ws.Cells.Item(12, 2).Value = "16:00:09"
$ws.Cells.Item(12, 3).Value = "16:00"
$ws.Cells.Item(12, 4).Value = "Living Room Main Door"
$ws.Cells.Item(12, 5).Value = "Image Captured"
$ws.Cells.Item(12, 6).Value = "Active"
